$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Remaining Space" column (D) for data rows 2-480 needs to be
# reduced by 5 (hotfix to account for electronic-equipment weight
# previously omitted from the calculation).
$lastRow = 480
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.Value = $cell.Value2 - 5
}
